$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 3.645393585217082

$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 10.29869402782916
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 659.8312095683168

$ws.Range("B4").Value = 0.01514828764759746
$ws.Range("C4").Value = 0.04240448674262143
$ws.Range("D4").Value = 689428.5527653177
$ws.Range("E4").Value = 616238.5361209477
$ws.Range("G4").Value = 1305667.14643904

$ws.Range("B5").Value = 0.6753301551942219
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 3.645393585217082

$ws.Range("B6").Value = 3.230985683306322
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 0.8054896365839992
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("G6").Value = 6.201049113329182

$ws.Range("B7").Value = 0.01514828764759746
$ws.Range("C7").Value = 0.3127903958511391
$ws.Range("D7").Value = 337.1190423067083
$ws.Range("E7").Value = 645.3272768299601
$ws.Range("G7").Value = 982.7742578201671

$ws.Range("B8").Value = 3.230985683306322
$ws.Range("C8").Value = 1.667794583268128
$ws.Range("D8").Value = 0.8054896365839992
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("G8").Value = 6.201049113329182

$ws.Range("B9").Value = 3.230985683306322
$ws.Range("C9").Value = 1.667794583268128
$ws.Range("D9").Value = 0.1575252929769615
$ws.Range("E9").Value = 0.496779210170732
$ws.Range("G9").Value = 5.553084769722144

$ws.Range("B10").Value = 0.6753301551942219
$ws.Range("C10").Value = 1.667794583268128
$ws.Range("D10").Value = 0.8054896365839992
$ws.Range("E10").Value = 0.496779210170732
$ws.Range("G10").Value = 3.645393585217082

$ws.Range("B11").Value = 0.3048080303191223
$ws.Range("C11").Value = 0.3127903958511391
$ws.Range("D11").Value = 0.1575252929769615
$ws.Range("E11").Value = 0.496779210170732
$ws.Range("G11").Value = 1.271902929317955

$ws.Range("B12").Value = 0.127881588408715
$ws.Range("C12").Value = 0.3127903958511391
$ws.Range("D12").Value = 0.8054896365839992
$ws.Range("E12").Value = 0.496779210170732
$ws.Range("G12").Value = 1.742940831014585

$ws.Range("B13").Value = 3.230985683306322
$ws.Range("C13").Value = 1.667794583268128
$ws.Range("D13").Value = 0.1575252929769615
$ws.Range("E13").Value = 0.496779210170732
$ws.Range("G13").Value = 5.553084769722144

$ws.Range("B14").Value = 0.00000001477430533469715
$ws.Range("C14").Value = 0.00007097389502863649
$ws.Range("D14").Value = 3.900430680208489
$ws.Range("E14").Value = 8.660232485948974
$ws.Range("G14").Value = 12.5607341548268

$ws.Range("B15").Value = 1.459612070389937
$ws.Range("C15").Value = 1.667794583268128
$ws.Range("D15").Value = 3.900430680208489
$ws.Range("E15").Value = 8.660232485948974
$ws.Range("G15").Value = 15.68806981981553

$ws.Range("B16").Value = 3.230985683306322
$ws.Range("C16").Value = 1.667794583268128
$ws.Range("D16").Value = 0.8054896365839992
$ws.Range("E16").Value = 8.660232485948974
$ws.Range("G16").Value = 14.36450238910742

$ws.Range("B17").Value = 3.230985683306322
$ws.Range("C17").Value = 0.3127903958511391
$ws.Range("D17").Value = 0.1575252929769615
$ws.Range("E17").Value = 0.496779210170732
$ws.Range("G17").Value = 4.198080582305154

$ws.Range("B18").Value = 3.230985683306322
$ws.Range("C18").Value = 10.29869402782916
$ws.Range("D18").Value = 337.1190423067083
$ws.Range("E18").Value = 8.660232485948974
$ws.Range("G18").Value = 359.3089545037927
